$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph and insert a new
# paragraph right after it containing the professor's name, styled as a
# bulleted list item (ListBullet).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r") -eq "Docente(s) Responsável(eis) ") {
        $target = $p
        break
    }
}

$target.Range.InsertParagraphAfter()
$newPara = $target.Next()
$newPara.Range.Text = "11079086 - Herlandí de Souza Andrade"
$newPara.Style = "ListBullet"
